# Update the value stored in A1 and move the active selection to D4,
# matching the target revision of Online.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1: 1 -> 2
$ws.Range("A1").Value = 2

# Active selection: A2 -> D4
$ws.Range("D4").Select()
